{"js": "// Update the division-problem cells in the table: each \"AAA\u00f7B=\" string is\n// replaced by a new \"CCC\u00f7D=\" string, matching the author's regenerated\n// worksheet numbers. Every search key below is unique in the document, so a\n// simple search-and-replace per pair is safe (no cross-hit ordering issues).\nconst replacements = [\n  [\"637\u00f75=\", \"343\u00f77=\"],\n  [\"735\u00f72=\", \"685\u00f76=\"],\n  [\"868\u00f72=\", \"985\u00f75=\"],\n  [\"577\u00f77=\", \"192\u00f72=\"],\n  [\"885\u00f73=\", \"995\u00f72=\"],\n  [\"638\u00f78=\", \"575\u00f72=\"],\n  [\"599\u00f79=\", \"404\u00f78=\"],\n  [\"796\u00f74=\", \"462\u00f78=\"],\n  [\"666\u00f73=\", \"997\u00f78=\"],\n  [\"566\u00f74=\", \"459\u00f75=\"],\n  [\"303\u00f72=\", \"658\u00f77=\"],\n  [\"879\u00f75=\", \"644\u00f78=\"],\n  [\"470\u00f77=\", \"545\u00f76=\"],\n  [\"159\u00f77=\", \"948\u00f72=\"],\n  [\"414\u00f77=\", \"172\u00f75=\"],\n  [\"144\u00f74=\", \"920\u00f77=\"],\n  [\"507\u00f77=\", \"391\u00f76=\"],\n  [\"248\u00f72=\", \"315\u00f74=\"],\n  [\"547\u00f77=\", \"533\u00f76=\"],\n  [\"370\u00f77=\", \"524\u00f72=\"],\n  [\"167\u00f78=\", \"925\u00f79=\"],\n  [\"516\u00f72=\", \"747\u00f72=\"],\n  [\"879\u00f76=\", \"790\u00f75=\"],\n  [\"877\u00f79=\", \"292\u00f72=\"],\n  [\"322\u00f79=\", \"101\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division-problem cells in the table: each \"AAA\u00f7B=\" string is\n# replaced by a new \"CCC\u00f7D=\" string, matching the author's regenerated\n# worksheet numbers. Every search key is unique in the document, so a plain\n# Find/Replace-All per pair is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"637\u00f75=\", \"343\u00f77=\"),\n    @(\"735\u00f72=\", \"685\u00f76=\"),\n    @(\"868\u00f72=\", \"985\u00f75=\"),\n    @(\"577\u00f77=\", \"192\u00f72=\"),\n    @(\"885\u00f73=\", \"995\u00f72=\"),\n    @(\"638\u00f78=\", \"575\u00f72=\"),\n    @(\"599\u00f79=\", \"404\u00f78=\"),\n    @(\"796\u00f74=\", \"462\u00f78=\"),\n    @(\"666\u00f73=\", \"997\u00f78=\"),\n    @(\"566\u00f74=\", \"459\u00f75=\"),\n    @(\"303\u00f72=\", \"658\u00f77=\"),\n    @(\"879\u00f75=\", \"644\u00f78=\"),\n    @(\"470\u00f77=\", \"545\u00f76=\"),\n    @(\"159\u00f77=\", \"948\u00f72=\"),\n    @(\"414\u00f77=\", \"172\u00f75=\"),\n    @(\"144\u00f74=\", \"920\u00f77=\"),\n    @(\"507\u00f77=\", \"391\u00f76=\"),\n    @(\"248\u00f72=\", \"315\u00f74=\"),\n    @(\"547\u00f77=\", \"533\u00f76=\"),\n    @(\"370\u00f77=\", \"524\u00f72=\"),\n    @(\"167\u00f78=\", \"925\u00f79=\"),\n    @(\"516\u00f72=\", \"747\u00f72=\"),\n    @(\"879\u00f76=\", \"790\u00f75=\"),\n    @(\"877\u00f79=\", \"292\u00f72=\"),\n    @(\"322\u00f79=\", \"101\u00f74=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n$d.Save()\n"}
